$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from an existing header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("C1").PasteSpecial(-4122)
$wsForecast.Range("D1").PasteSpecial(-4122)

# Data rows
$data = @(
    @(45207.99999999999, 16, -11.50561084966274,  45.89146067486202),
    @(45221.99999999999, 17, -13.48330411786288,  45.43518120984754),
    @(45333.99999999999, 17, -10.19724211774583,  47.61170462966027),
    @(45592.99999999999, 19, -11.66857187981909,  49.20623876276937),
    @(45599.99999999999, 19, -8.464684533340318,  49.60722643987874),
    @(45606.99999999999, 19, -9.2308282544021,    50.13109318780788),
    @(45613.99999999999, 19, -9.797072872490753,  50.97558154074635),
    @(45620.99999999999, 19, -11.0085822869318,   49.98985627461941),
    @(45627.99999999999, 20, -8.465564897690147,  46.31093521725408),
    @(45634.99999999999, 20, -10.36173839794875,  48.80614528936207),
    @(45641.99999999999, 20, -8.096390144988911,  49.59335048026462),
    @(45648.99999999999, 20, -8.535642537662014,  47.44137533952357)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Copy the date-style formatting (s="2") from the "Order Week" date column onto column A data rows
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()
